$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at row 41, shifting existing rows 41-54 down to 42-55
$ws.Rows.Item(41).Insert()

# Populate the new row 41 with the new parameter entry
# (order matches shared-string creation order of the target file: A, D, E, H, F, G)
$ws.Range("A41").Value = "reconstruction_functional_network.reconstructionMethod"
$ws.Range("D41").Value = "reconstruction_functional_network"
$ws.Range("E41").Value = "char"
$ws.Range("H41").Value = "Functional connectivity estimation method used."
$ws.Range("F41").Value = "isfunction nonempty"
$ws.Range("G41").Value = "standard"

# Match style used by other rows in F/G columns (numFmtId 49 - text format, "@")
$ws.Range("F41").NumberFormat = "@"
$ws.Range("G41").NumberFormat = "@"

# Update selection like in the target file
$ws.Range("F41").Select()
